$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder and update existing rows (2-4), then add new rows (5-7)
# Target layout:
# Row2: Wrong_Entity_NonEvent_as_Event = 76
# Row3: Correct = 38
# Row4: Wrong_Entity_Event_as_NonEvent = 35
# Row5: Wrong_Tag_I_as_B = 1
# Row6: Wrong_Tag_B_as_I = 1
# Row7: Wrong_Tag_S_as_B = 1

$ws.Range("A2").Value = "Wrong_Entity_NonEvent_as_Event"
$ws.Range("B2").Value = 76

$ws.Range("A3").Value = "Correct"
$ws.Range("B3").Value = 38

$ws.Range("A4").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B4").Value = 35

$ws.Range("A5").Value = "Wrong_Tag_I_as_B"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "Wrong_Tag_B_as_I"
$ws.Range("B6").Value = 1

$ws.Range("A7").Value = "Wrong_Tag_S_as_B"
$ws.Range("B7").Value = 1

# Apply the same style as other A-column cells (A2:A4) to the new rows A5:A7
$ws.Range("A4").Copy()
$ws.Range("A5:A7").PasteSpecial(-4122)
